$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data update from DGS's 2021/12/06, 2021/12/07 and 2021/12/08 reports:
# two new rows are appended right after the current last row (115), for
# report dates 2021/12/06 (row 116) and 2021/12/08 (row 117).

# First, replicate the formatting/style of the last existing row onto the
# two new rows so the new cells pick up the same styles (s="1" for the date
# column, s="2" for the numeric columns) as the rest of the table.
$ws.Range("A115:E115").Copy()
$ws.Range("A116:E117").PasteSpecial(-4122)

# Column A stores the report date as text (a shared string formatted to
# look like yyyy/mm/dd), not a real date serial number. Setting the value
# directly on a date-formatted cell would make Excel auto-convert it into a
# date serial, so instead set it on a plain ("@"/Text) formatted cell and
# copy the resulting text value across - this keeps the cell's style (and
# number format) untouched while avoiding the automatic date conversion.
$ws.Range("A116").NumberFormat = "@"
$ws.Range("A116").Value = "2021/12/06"
$ws.Range("A116").NumberFormat = "yyyy/mm/dd"

$ws.Range("A117").NumberFormat = "@"
$ws.Range("A117").Value = "2021/12/08"
$ws.Range("A117").NumberFormat = "yyyy/mm/dd"

$ws.Range("B116").Value = 410.4
$ws.Range("C116").Value = 413.9
$ws.Range("D116").Value = 1.1
$ws.Range("E116").Value = 1.11

$ws.Range("B117").Value = 438.4
$ws.Range("C117").Value = 442.1
$ws.Range("D117").Value = 1.11
$ws.Range("E117").Value = 1.11

# Move the active cell/selection to the next empty row, same as the
# original workbook kept the selection on the row right after the data.
$ws.Range("A118").Select()
